# Insert a new row at position 70 (shifting existing rows 70-149 down to 71-150)
# and populate the new row 70 with the latest price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("70:70").Insert()

$ws.Range("A70").Value = 5
$ws.Range("B70").Value = "Macroferia Regional de Talca"
$ws.Range("C70").Value = "Maule"
$ws.Range("D70").Value = 44994
$ws.Range("E70").Value = 7
$ws.Range("F70").Value = 100112001
$ws.Range("G70").Value = "Berenjena"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 8000
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = 8000
$ws.Range("N70").Value = "$/caja 50 unidades"
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 160
$ws.Range("Q70").Value = 50
$ws.Range("R70").Value = "Hortaliza"
